$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.236.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.564.94'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.93'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.490'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.18'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.788.07'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.571.52'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.76'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.260.14'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.99'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.42'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.14'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.29'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.80'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.107'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.74%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.02'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.60%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.455.28'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.12%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +5.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.62'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.34'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.540'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.815'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.34'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.40'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.74'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.703.96'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.77'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.29%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0947'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.33%  '
